# Scheduled market-data refresh (Universalis price pull) for the Leve-profit workbook.
# Recomputed currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) per job sheet
# to reflect the latest Eorzea market snapshot. No structural/layout changes.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1480.2307
$ws.Range("I80").Value = 649.5
$ws.Range("J80").Value = 1849.4445
$ws.Range("K80").Value = 1948.5
$ws.Range("L80").Value = 5548.333500000001
$ws.Range("M80").Value = -950.5
$ws.Range("N80").Value = -7544.333500000001
$ws.Range("H83").Value = 1480.2307
$ws.Range("I83").Value = 649.5
$ws.Range("J83").Value = 1849.4445
$ws.Range("K83").Value = 5845.5
$ws.Range("L83").Value = 16645.0005
$ws.Range("M83").Value = -853.5
$ws.Range("N83").Value = -26629.0005
$ws.Range("H86").Value = 8933.333000000001
$ws.Range("J86").Value = 9333.333000000001
$ws.Range("L86").Value = 9333.333000000001
$ws.Range("N86").Value = -11579.333
$ws.Range("H88").Value = 2444.2222
$ws.Range("I88").Value = 2810.5
$ws.Range("J88").Value = 2339.5715
$ws.Range("K88").Value = 2810.5
$ws.Range("L88").Value = 2339.5715
$ws.Range("M88").Value = -2404.5
$ws.Range("N88").Value = -3151.5715
$ws.Range("H89").Value = 8933.333000000001
$ws.Range("J89").Value = 9333.333000000001
$ws.Range("L89").Value = 46666.665
$ws.Range("N89").Value = -57898.665
$ws.Range("H91").Value = 2444.2222
$ws.Range("I91").Value = 2810.5
$ws.Range("J91").Value = 2339.5715
$ws.Range("K91").Value = 2810.5
$ws.Range("L91").Value = 2339.5715
$ws.Range("M91").Value = -1406.5
$ws.Range("N91").Value = -5147.5715
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H138").Value = 2398.1155
$ws.Range("I138").Value = 1491.7222
$ws.Range("J138").Value = 4437.5
$ws.Range("K138").Value = 4475.1666
$ws.Range("L138").Value = 13312.5
$ws.Range("M138").Value = 664.8334000000004
$ws.Range("N138").Value = -23592.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4478.077
$ws.Range("I2").Value = 1542.2
$ws.Range("J2").Value = 6313
$ws.Range("K2").Value = 1542.2
$ws.Range("L2").Value = 6313
$ws.Range("M2").Value = -1429.2
$ws.Range("N2").Value = -6539
$ws.Range("H74").Value = 1426.8
$ws.Range("I74").Value = 1426.8
$ws.Range("K74").Value = 1426.8
$ws.Range("M74").Value = -552.8
$ws.Range("H77").Value = 1426.8
$ws.Range("I77").Value = 1426.8
$ws.Range("K77").Value = 7134
$ws.Range("M77").Value = -2766
$ws.Range("H88").Value = 2544.5
$ws.Range("J88").Value = 3105.8235
$ws.Range("L88").Value = 3105.8235
$ws.Range("N88").Value = -3917.8235
$ws.Range("H91").Value = 2544.5
$ws.Range("J91").Value = 3105.8235
$ws.Range("L91").Value = 3105.8235
$ws.Range("N91").Value = -5913.8235
$ws.Range("H116").Value = 4478.077
$ws.Range("I116").Value = 1542.2
$ws.Range("J116").Value = 6313
$ws.Range("K116").Value = 1542.2
$ws.Range("L116").Value = 6313
$ws.Range("M116").Value = 751.8
$ws.Range("N116").Value = -10901
$ws.Range("H122").Value = 2946.0908
$ws.Range("J122").Value = 4019.25
$ws.Range("L122").Value = 12057.75
$ws.Range("N122").Value = -16957.75

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4478.077
$ws.Range("I3").Value = 1542.2
$ws.Range("J3").Value = 6313
$ws.Range("K3").Value = 1542.2
$ws.Range("L3").Value = 6313
$ws.Range("M3").Value = -1428.2
$ws.Range("N3").Value = -6541
$ws.Range("H133").Value = 70709
$ws.Range("I133").Value = 70709
$ws.Range("K133").Value = 70709
$ws.Range("M133").Value = -65649

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 8316.75
$ws.Range("I99").Value = 8163.8184
$ws.Range("K99").Value = 8163.8184
$ws.Range("M99").Value = -6665.8184
$ws.Range("H126").Value = 8316.75
$ws.Range("I126").Value = 8163.8184
$ws.Range("K126").Value = 24491.4552
$ws.Range("M126").Value = -22021.4552

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 119
$ws.Range("J12").Value = 148
$ws.Range("L12").Value = 444
$ws.Range("N12").Value = -790
$ws.Range("H40").Value = 118.44444
$ws.Range("I40").Value = 128.25
$ws.Range("J40").Value = 40
$ws.Range("K40").Value = 513
$ws.Range("L40").Value = 160
$ws.Range("M40").Value = -444
$ws.Range("N40").Value = -298

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2499.125
$ws.Range("I80").Value = 2513.2856
$ws.Range("J80").Value = 2400
$ws.Range("K80").Value = 2513.2856
$ws.Range("L80").Value = 2400
$ws.Range("M80").Value = -1515.2856
$ws.Range("N80").Value = -4396
$ws.Range("H83").Value = 2499.125
$ws.Range("I83").Value = 2513.2856
$ws.Range("J83").Value = 2400
$ws.Range("K83").Value = 12566.428
$ws.Range("L83").Value = 12000
$ws.Range("M83").Value = -7574.428
$ws.Range("N83").Value = -21984

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2166.6667
$ws.Range("I46").Value = 1642.8572
$ws.Range("J46").Value = 2625
$ws.Range("K46").Value = 1642.8572
$ws.Range("L46").Value = 2625
$ws.Range("M46").Value = -1454.8572
$ws.Range("N46").Value = -3001
$ws.Range("H82").Value = 3000
$ws.Range("I82").Value = 3000
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 3000
$ws.Range("L82").Value = 3000
$ws.Range("M82").Value = -2639
$ws.Range("N82").Value = -3722
$ws.Range("H85").Value = 3000
$ws.Range("I85").Value = 3000
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 3000
$ws.Range("L85").Value = 3000
$ws.Range("M85").Value = -1752
$ws.Range("N85").Value = -5496

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3599.8
$ws.Range("I81").Value = 3999.75
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 7999.5
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -6938.5
$ws.Range("N81").Value = -6122
$ws.Range("H84").Value = 3599.8
$ws.Range("I84").Value = 3999.75
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 39997.5
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -34693.5
$ws.Range("N84").Value = -30608
$ws.Range("H126").Value = 2448.5
$ws.Range("I126").Value = 2389.2727
$ws.Range("K126").Value = 7167.8181
$ws.Range("M126").Value = -4697.8181
